$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 742.6042
$ws.Range("J129").Value = 972
$ws.Range("L129").Value = 2916
$ws.Range("N129").Value = -12916

$ws.Range("H137").Value = 33425.03
$ws.Range("I137").Value = 46446.727
$ws.Range("K137").Value = 139340.181
$ws.Range("M137").Value = -136790.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 37637.332
$ws.Range("J7").Value = 37637.332
$ws.Range("L7").Value = 37637.332
$ws.Range("N7").Value = -37865.332

$ws.Range("H32").Value = 5428.3696
$ws.Range("I32").Value = 5545.1025
$ws.Range("J32").Value = 4778
$ws.Range("K32").Value = 5545.1025
$ws.Range("L32").Value = 4778
$ws.Range("M32").Value = -5258.1025
$ws.Range("N32").Value = -5352

$ws.Range("H74").Value = 19046.357
$ws.Range("I74").Value = 21894.209
$ws.Range("J74").Value = 1959.25
$ws.Range("K74").Value = 21894.209
$ws.Range("L74").Value = 1959.25
$ws.Range("M74").Value = -21020.209
$ws.Range("N74").Value = -3707.25

$ws.Range("H77").Value = 19046.357
$ws.Range("I77").Value = 21894.209
$ws.Range("J77").Value = 1959.25
$ws.Range("K77").Value = 109471.045
$ws.Range("L77").Value = 9796.25
$ws.Range("M77").Value = -105103.045
$ws.Range("N77").Value = -18532.25

$ws.Range("H122").Value = 1013.4138
$ws.Range("I122").Value = 902.5417
$ws.Range("J122").Value = 1545.6
$ws.Range("K122").Value = 2707.6251
$ws.Range("L122").Value = 4636.799999999999
$ws.Range("M122").Value = -257.6251000000002
$ws.Range("N122").Value = -9536.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1828.7142
$ws.Range("I94").Value = 590
$ws.Range("J94").Value = 3067.4285
$ws.Range("K94").Value = 590
$ws.Range("L94").Value = 3067.4285
$ws.Range("M94").Value = -139
$ws.Range("N94").Value = -3969.4285

$ws.Range("H99").Value = 938.4211
$ws.Range("I99").Value = 1002.5
$ws.Range("J99").Value = 828.5714
$ws.Range("K99").Value = 1002.5
$ws.Range("L99").Value = 828.5714
$ws.Range("M99").Value = 495.5
$ws.Range("N99").Value = -3824.5714

$ws.Range("H134").Value = 17313.883
$ws.Range("I134").Value = 932.4
$ws.Range("J134").Value = 62818
$ws.Range("K134").Value = 2797.2
$ws.Range("L134").Value = 188454
$ws.Range("M134").Value = -262.1999999999998
$ws.Range("N134").Value = -193524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6922.768
$ws.Range("I31").Value = 4929.5166
$ws.Range("K31").Value = 4929.5166
$ws.Range("M31").Value = -4634.5166

$ws.Range("H34").Value = 6922.768
$ws.Range("I34").Value = 4929.5166
$ws.Range("K34").Value = 4929.5166
$ws.Range("M34").Value = -4727.5166

$ws.Range("H58").Value = 901.5507
$ws.Range("I58").Value = 667.9804
$ws.Range("J58").Value = 1563.3334
$ws.Range("K58").Value = 667.9804
$ws.Range("L58").Value = 1563.3334
$ws.Range("M58").Value = -464.9804
$ws.Range("N58").Value = -1969.3334

$ws.Range("H132").Value = 837.2449
$ws.Range("I132").Value = 888.9773
$ws.Range("J132").Value = 382
$ws.Range("K132").Value = 2666.9319
$ws.Range("L132").Value = 1146
$ws.Range("M132").Value = -136.9319
$ws.Range("N132").Value = -6206

$ws.Range("H136").Value = 901.5507
$ws.Range("I136").Value = 667.9804
$ws.Range("J136").Value = 1563.3334
$ws.Range("K136").Value = 2003.9412
$ws.Range("L136").Value = 4690.0002
$ws.Range("M136").Value = 546.0587999999998
$ws.Range("N136").Value = -9790.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 81.935486
$ws.Range("I12").Value = 233.14285
$ws.Range("J12").Value = 37.833332
$ws.Range("K12").Value = 699.4285500000001
$ws.Range("L12").Value = 113.499996
$ws.Range("M12").Value = -526.4285500000001
$ws.Range("N12").Value = -459.499996

$ws.Range("H64").Value = 1517655.6
$ws.Range("I64").Value = 1169
$ws.Range("J64").Value = 2167578.5
$ws.Range("K64").Value = 3507
$ws.Range("L64").Value = 6502735.5
$ws.Range("M64").Value = -3237
$ws.Range("N64").Value = -6503275.5

$ws.Range("H67").Value = 1517655.6
$ws.Range("I67").Value = 1169
$ws.Range("J67").Value = 2167578.5
$ws.Range("K67").Value = 3507
$ws.Range("L67").Value = 6502735.5
$ws.Range("M67").Value = -2571
$ws.Range("N67").Value = -6504607.5

$ws.Range("H120").Value = 13558.889
$ws.Range("I120").Value = 4257.5
$ws.Range("J120").Value = 15176.521
$ws.Range("K120").Value = 12772.5
$ws.Range("L120").Value = 45529.563
$ws.Range("M120").Value = -7934.5
$ws.Range("N120").Value = -55205.563

$ws.Range("H122").Value = 331.34885
$ws.Range("I122").Value = 254.72
$ws.Range("J122").Value = 437.77777
$ws.Range("K122").Value = 2292.48
$ws.Range("L122").Value = 3939.99993
$ws.Range("M122").Value = 157.52
$ws.Range("N122").Value = -8839.99993

$ws.Range("H131").Value = 73530030
$ws.Range("J131").Value = 125000776
$ws.Range("L131").Value = 375002328
$ws.Range("N131").Value = -375012408

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 27333.54
$ws.Range("I132").Value = 1901.4445
$ws.Range("J132").Value = 49132.477
$ws.Range("K132").Value = 5704.333500000001
$ws.Range("L132").Value = 147397.431
$ws.Range("M132").Value = -3174.333500000001
$ws.Range("N132").Value = -152457.431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5004.1763
$ws.Range("I7").Value = 4210.4287
$ws.Range("J7").Value = 5559.8
$ws.Range("K7").Value = 4210.4287
$ws.Range("L7").Value = 5559.8
$ws.Range("M7").Value = -4098.4287
$ws.Range("N7").Value = -5783.8

$ws.Range("H126").Value = 5004.1763
$ws.Range("I126").Value = 4210.4287
$ws.Range("J126").Value = 5559.8
$ws.Range("K126").Value = 12631.2861
$ws.Range("L126").Value = 16679.4
$ws.Range("M126").Value = -10161.2861
$ws.Range("N126").Value = -21619.4

$ws.Range("H136").Value = 176587.86
$ws.Range("I136").Value = 244670.56
$ws.Range("J136").Value = 2125.9375
$ws.Range("K136").Value = 734011.6799999999
$ws.Range("L136").Value = 6377.8125
$ws.Range("M136").Value = -731461.6799999999
$ws.Range("N136").Value = -11477.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 731.6799999999999
$ws.Range("I126").Value = 604.1905
$ws.Range("K126").Value = 1812.5715
$ws.Range("M126").Value = 657.4285

$ws.Range("H132").Value = 1712.1794
$ws.Range("I132").Value = 306.87756
$ws.Range("J132").Value = 4086.6553
$ws.Range("K132").Value = 920.6326800000001
$ws.Range("L132").Value = 12259.9659
$ws.Range("M132").Value = 1609.36732
$ws.Range("N132").Value = -17319.9659

$ws.Range("H136").Value = 895340.9
$ws.Range("I136").Value = 1190916.4
$ws.Range("J136").Value = 323259.34
$ws.Range("K136").Value = 3572749.2
$ws.Range("L136").Value = 969778.02
$ws.Range("M136").Value = -3570199.2
$ws.Range("N136").Value = -974878.02
